$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 44.855544
$ws.Cells.Item(2, 8).Value = 134.566632
$ws.Cells.Item(2, 9).Value = 0.9269627513664965
$ws.Cells.Item(2, 10).Value = 0.9269627513664968
$ws.Cells.Item(2, 13).Value = 70.23436
$ws.Cells.Item(2, 14).Value = 210.70308
$ws.Cells.Item(2, 15).Value = 0.7023186840741513
$ws.Cells.Item(2, 16).Value = 0.7023186840741513
$ws.Cells.Item(2, 17).Value = 3150.40042529184
$ws.Cells.Item(2, 18).Value = 28353.60382762656
$ws.Cells.Item(2, 19).Value = 0.6510232597254726
$ws.Cells.Item(2, 20).Value = 0.6510232597254727
$ws.Cells.Item(3, 7).Value = 44.855544
$ws.Cells.Item(3, 8).Value = 134.566632
$ws.Cells.Item(3, 9).Value = 0.9269627513664965
$ws.Cells.Item(3, 10).Value = 0.9269627513664968
$ws.Cells.Item(3, 15).Value = 0.06551129587759326
$ws.Cells.Item(3, 16).Value = 0.06551129587759325
$ws.Cells.Item(3, 17).Value = 293.864906450928
$ws.Cells.Item(3, 18).Value = 2644.784158058352
$ws.Cells.Item(3, 19).Value = 0.06072653107227847
$ws.Cells.Item(3, 20).Value = 0.06072653107227847
$ws.Cells.Item(4, 7).Value = 44.855544
$ws.Cells.Item(4, 8).Value = 134.566632
$ws.Cells.Item(4, 9).Value = 0.9269627513664965
$ws.Cells.Item(4, 10).Value = 0.9269627513664968
$ws.Cells.Item(4, 13).Value = 4.268944666666666
$ws.Cells.Item(4, 14).Value = 12.806834
$ws.Cells.Item(4, 15).Value = 0.04268793224112385
$ws.Cells.Item(4, 16).Value = 0.04268793224112385
$ws.Cells.Item(4, 17).Value = 191.485835329232
$ws.Cells.Item(4, 18).Value = 1723.372517963088
$ws.Cells.Item(4, 19).Value = 0.03957012312037874
$ws.Cells.Item(4, 20).Value = 0.03957012312037875
$ws.Cells.Item(5, 7).Value = 44.855544
$ws.Cells.Item(5, 8).Value = 134.566632
$ws.Cells.Item(5, 9).Value = 0.9269627513664965
$ws.Cells.Item(5, 10).Value = 0.9269627513664968
$ws.Cells.Item(5, 13).Value = 18.948881
$ws.Cells.Item(5, 14).Value = 56.846643
$ws.Cells.Item(5, 15).Value = 0.1894820878071316
$ws.Cells.Item(5, 16).Value = 0.1894820878071315
$ws.Cells.Item(5, 17).Value = 849.962365446264
$ws.Cells.Item(5, 18).Value = 7649.661289016376
$ws.Cells.Item(5, 19).Value = 0.1756428374483668
$ws.Cells.Item(5, 20).Value = 0.1756428374483668
$ws.Cells.Item(6, 9).Value = 0.03026428998407557
$ws.Cells.Item(6, 10).Value = 0.03026428998407558
$ws.Cells.Item(6, 13).Value = 70.23436
$ws.Cells.Item(6, 14).Value = 210.70308
$ws.Cells.Item(6, 15).Value = 0.7023186840741513
$ws.Cells.Item(6, 16).Value = 0.7023186840741513
$ws.Cells.Item(6, 17).Value = 102.85702623588
$ws.Cells.Item(6, 18).Value = 925.7132361229201
$ws.Cells.Item(6, 19).Value = 0.02125517631605447
$ws.Cells.Item(6, 20).Value = 0.02125517631605448
$ws.Cells.Item(7, 9).Value = 0.03026428998407557
$ws.Cells.Item(7, 10).Value = 0.03026428998407558
$ws.Cells.Item(7, 15).Value = 0.06551129587759326
$ws.Cells.Item(7, 16).Value = 0.06551129587759325
$ws.Cells.Item(7, 19).Value = 0.001982652855672057
$ws.Cells.Item(7, 20).Value = 0.001982652855672057
$ws.Cells.Item(8, 9).Value = 0.03026428998407557
$ws.Cells.Item(8, 10).Value = 0.03026428998407558
$ws.Cells.Item(8, 13).Value = 4.268944666666666
$ws.Cells.Item(8, 14).Value = 12.806834
$ws.Cells.Item(8, 15).Value = 0.04268793224112385
$ws.Cells.Item(8, 16).Value = 0.04268793224112385
$ws.Cells.Item(8, 17).Value = 6.251796892274
$ws.Cells.Item(8, 18).Value = 56.266172030466
$ws.Cells.Item(8, 19).Value = 0.001291919960165941
$ws.Cells.Item(8, 20).Value = 0.001291919960165942
$ws.Cells.Item(9, 9).Value = 0.03026428998407557
$ws.Cells.Item(9, 10).Value = 0.03026428998407558
$ws.Cells.Item(9, 13).Value = 18.948881
$ws.Cells.Item(9, 14).Value = 56.846643
$ws.Cells.Item(9, 15).Value = 0.1894820878071316
$ws.Cells.Item(9, 16).Value = 0.1894820878071315
$ws.Cells.Item(9, 17).Value = 27.75031409352301
$ws.Cells.Item(9, 18).Value = 249.752826841707
$ws.Cells.Item(9, 19).Value = 0.0057345408521831
$ws.Cells.Item(9, 20).Value = 0.005734540852183101
$ws.Cells.Item(10, 7).Value = 1.967437666666666
$ws.Cells.Item(10, 8).Value = 5.902312999999999
$ws.Cells.Item(10, 9).Value = 0.04065810533109158
$ws.Cells.Item(10, 10).Value = 0.0406581053310916
$ws.Cells.Item(10, 13).Value = 70.23436
$ws.Cells.Item(10, 14).Value = 210.70308
$ws.Cells.Item(10, 15).Value = 0.7023186840741513
$ws.Cells.Item(10, 16).Value = 0.7023186840741513
$ws.Cells.Item(10, 17).Value = 138.1817253582266
$ws.Cells.Item(10, 18).Value = 1243.63552822404
$ws.Cells.Item(10, 19).Value = 0.02855494703308048
$ws.Cells.Item(10, 20).Value = 0.02855494703308049
$ws.Cells.Item(11, 7).Value = 1.967437666666666
$ws.Cells.Item(11, 8).Value = 5.902312999999999
$ws.Cells.Item(11, 9).Value = 0.04065810533109158
$ws.Cells.Item(11, 10).Value = 0.0406581053310916
$ws.Cells.Item(11, 15).Value = 0.06551129587759326
$ws.Cells.Item(11, 16).Value = 0.06551129587759325
$ws.Cells.Item(11, 17).Value = 12.88939636676866
$ws.Cells.Item(11, 18).Value = 116.004567300918
$ws.Cells.Item(11, 19).Value = 0.002663565168167493
$ws.Cells.Item(11, 20).Value = 0.002663565168167493
$ws.Cells.Item(12, 7).Value = 1.967437666666666
$ws.Cells.Item(12, 8).Value = 5.902312999999999
$ws.Cells.Item(12, 9).Value = 0.04065810533109158
$ws.Cells.Item(12, 10).Value = 0.0406581053310916
$ws.Cells.Item(12, 13).Value = 4.268944666666666
$ws.Cells.Item(12, 14).Value = 12.806834
$ws.Cells.Item(12, 15).Value = 0.04268793224112385
$ws.Cells.Item(12, 16).Value = 0.04268793224112385
$ws.Cells.Item(12, 17).Value = 8.398882534115776
$ws.Cells.Item(12, 18).Value = 75.58994280704198
$ws.Cells.Item(12, 19).Value = 0.001735610445426114
$ws.Cells.Item(12, 20).Value = 0.001735610445426115
$ws.Cells.Item(13, 7).Value = 1.967437666666666
$ws.Cells.Item(13, 8).Value = 5.902312999999999
$ws.Cells.Item(13, 9).Value = 0.04065810533109158
$ws.Cells.Item(13, 10).Value = 0.0406581053310916
$ws.Cells.Item(13, 13).Value = 18.948881
$ws.Cells.Item(13, 14).Value = 56.846643
$ws.Cells.Item(13, 15).Value = 0.1894820878071316
$ws.Cells.Item(13, 16).Value = 0.1894820878071315
$ws.Cells.Item(13, 17).Value = 37.28074222058433
$ws.Cells.Item(13, 18).Value = 335.526679985259
$ws.Cells.Item(13, 19).Value = 0.0077039826844175
$ws.Cells.Item(13, 20).Value = 0.007703982684417502
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.1023373333333333
$ws.Cells.Item(14, 8).Value = 0.307012
$ws.Cells.Item(14, 9).Value = 0.002114853318336234
$ws.Cells.Item(14, 10).Value = 0.002114853318336234
$ws.Cells.Item(14, 13).Value = 70.23436
$ws.Cells.Item(14, 14).Value = 210.70308
$ws.Cells.Item(14, 15).Value = 0.7023186840741513
$ws.Cells.Item(14, 16).Value = 0.7023186840741513
$ws.Cells.Item(14, 17).Value = 7.187597110773333
$ws.Cells.Item(14, 18).Value = 64.68837399696
$ws.Cells.Item(14, 19).Value = 0.001485300999543756
$ws.Cells.Item(14, 20).Value = 0.001485300999543756
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.1023373333333333
$ws.Cells.Item(15, 8).Value = 0.307012
$ws.Cells.Item(15, 9).Value = 0.002114853318336234
$ws.Cells.Item(15, 10).Value = 0.002114853318336234
$ws.Cells.Item(15, 15).Value = 0.06551129587759326
$ws.Cells.Item(15, 16).Value = 0.06551129587759325
$ws.Cells.Item(15, 17).Value = 0.6704489167813333
$ws.Cells.Item(15, 18).Value = 6.034040251032
$ws.Cells.Item(15, 19).Value = 0.0001385467814752349
$ws.Cells.Item(15, 20).Value = 0.0001385467814752349
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.1023373333333333
$ws.Cells.Item(16, 8).Value = 0.307012
$ws.Cells.Item(16, 9).Value = 0.002114853318336234
$ws.Cells.Item(16, 10).Value = 0.002114853318336234
$ws.Cells.Item(16, 13).Value = 4.268944666666666
$ws.Cells.Item(16, 14).Value = 12.806834
$ws.Cells.Item(16, 15).Value = 0.04268793224112385
$ws.Cells.Item(16, 16).Value = 0.04268793224112385
$ws.Cells.Item(16, 17).Value = 0.4368724133342222
$ws.Cells.Item(16, 18).Value = 3.931851720007999
$ws.Cells.Item(16, 19).Value = 0.00009027871515305309
$ws.Cells.Item(16, 20).Value = 0.00009027871515305312
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.1023373333333333
$ws.Cells.Item(17, 8).Value = 0.307012
$ws.Cells.Item(17, 9).Value = 0.002114853318336234
$ws.Cells.Item(17, 10).Value = 0.002114853318336234
$ws.Cells.Item(17, 13).Value = 18.948881
$ws.Cells.Item(17, 14).Value = 56.846643
$ws.Cells.Item(17, 15).Value = 0.1894820878071316
$ws.Cells.Item(17, 16).Value = 0.1894820878071315
$ws.Cells.Item(17, 17).Value = 1.939177951190667
$ws.Cells.Item(17, 18).Value = 17.452601560716
$ws.Cells.Item(17, 19).Value = 0.0004007268221641899
$ws.Cells.Item(17, 20).Value = 0.0004007268221641899
